$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a plain number (e.g. "331.42") need the
# number format forced to Text first, otherwise Excel auto-converts the
# entry to a real number (losing the original "27.463.51"-style string
# formatting / precision). ClearFormats() afterwards drops the temporary
# "@" number format so the cell keeps the sheets default (unstyled) look.

$ws.Range('D2').Value = '27.472.20'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '1.832.63'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  -0.75%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '331.42'
$r.ClearFormats()
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('E6').Value = '  -0.76%  '
$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '0.4621'
$r.ClearFormats()
$ws.Range('E7').Value = '  -2.96%  '
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.3836'
$r.ClearFormats()
$ws.Range('E8').Value = '  -2.34%  '
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '46.67'
$r.ClearFormats()
$ws.Range('E9').Value = '  -0.56%  '
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.07891'
$r.ClearFormats()
$ws.Range('E10').Value = '  -1.36%  '
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '0.9725'
$r.ClearFormats()
$ws.Range('E11').Value = '  -3.81%  '
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '21.13'
$r.ClearFormats()
$ws.Range('E12').Value = '  -2.86%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.842.03'
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '5.889'
$r.ClearFormats()
$ws.Range('E14').Value = '  -2.12%  '
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '7.063'
$r.ClearFormats()
$ws.Range('E15').Value = '  -1.47%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '1.002'
$r.ClearFormats()
$ws.Range('E16').Value = '  -0.95%  '
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '88.05'
$r.ClearFormats()
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('E18').Value = '  -1.43%  '
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '0.00001030'
$r.ClearFormats()
$ws.Range('E19').Value = '  -1.57%  '
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '17.03'
$r.ClearFormats()
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('D22').Value = '27.475.11'
$ws.Range('E22').Value = '  -1.44%  '
$ws.Range('E23').Value = '  -2.59%  '
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '10.82'
$r.ClearFormats()
$ws.Range('E24').Value = '  -1.22%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '2.310'
$r.ClearFormats()
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.054.96'
$ws.Range('E26').Value = '  -2.14%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '157.21'
$r.ClearFormats()
$ws.Range('E27').Value = '  -0.55%  '
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '19.39'
$r.ClearFormats()
$ws.Range('E28').Value = '  -1.82%  '
$ws.Range('E29').Value = '  -1.45%  '
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '5.305'
$r.ClearFormats()
$ws.Range('E30').Value = '  -2.63%  '
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '119.12'
$r.ClearFormats()
$ws.Range('E31').Value = '  -1.88%  '
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '0.9577'
$r.ClearFormats()
$ws.Range('E32').Value = '  -1.36%  '
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '0.09293'
$r.ClearFormats()
$ws.Range('E33').Value = '  -2.09%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '3.567'
$r.ClearFormats()
$ws.Range('E34').Value = '  -1.82%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '5.243'
$r.ClearFormats()
$ws.Range('E35').Value = '  -1.52%  '
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '1.316'
$r.ClearFormats()
$ws.Range('E36').Value = '  -2.45%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '0.05940'
$r.ClearFormats()
$ws.Range('E37').Value = '  -2.11%  '
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '0.02198'
$r.ClearFormats()
$ws.Range('E38').Value = '  -1.40%  '
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '8.079'
$r.ClearFormats()
$ws.Range('E39').Value = '  -0.83%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '1.161'
$r.ClearFormats()
$ws.Range('E40').Value = '  -3.74%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.5806'
$r.ClearFormats()
$ws.Range('E41').Value = '  -2.43%  '
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '0.1841'
$r.ClearFormats()
$ws.Range('E42').Value = '  -2.64%  '
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '10.03'
$r.ClearFormats()
$ws.Range('E43').Value = '  -2.67%  '
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '1.286'
$r.ClearFormats()
$ws.Range('E44').Value = '  +2.41%  '
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '0.5491'
$r.ClearFormats()
$ws.Range('E45').Value = '  -2.62%  '
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '11.98'
$r.ClearFormats()
$ws.Range('E46').Value = '  -0.98%  '
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '1.874'
$r.ClearFormats()
$ws.Range('E47').Value = '  -2.37%  '
$ws.Range('E48').Value = '  -1.89%  '
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '110.51'
$r.ClearFormats()
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('E50').Value = '  -2.49%  '
$ws.Range('E51').Value = '  -0.83%  '
